$wb = $excel.ActiveWorkbook

# Update OFF sheet (row 2: Short Att, Short Comp, Deep Att, Deep Comp)
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 347
$wsOff.Range("C2").Value = 269
$wsOff.Range("D2").Value = 90
$wsOff.Range("E2").Value = 46

# Update DEF sheet (row 2: Short Att, Short Comp, Deep Att, Deep Comp, Short Int)
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 487
$wsDef.Range("C2").Value = 340
$wsDef.Range("D2").Value = 115
$wsDef.Range("E2").Value = 51
$wsDef.Range("F2").Value = 7
